# #272 Ajout d'un scenario de recherche de l'offre d'un professionnel avec un ID Nat PS
# - Bump the "Date" metadata value on the Metadata sheet.
# - Swap the two "Mapping" columns (AK / AL) on the Elements sheet: the
#   "Mapping: Spécification métier vers l'extension ROR NbTemporarySocialHelpPlace"
#   column now comes before the "Mapping: RIM Mapping" column, so both the header
#   text and the per-row data (and the column widths that go with them) move
#   from column AK to AL and vice versa.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) --------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 6
$colLeft = 37  # AK - currently "Mapping: RIM Mapping"
$colRight = 38 # AL - currently "Mapping: Spécification métier..."

for ($r = 1; $r -le $lastRow; $r++) {
    $leftCell = $elements.Cells.Item($r, $colLeft)
    $rightCell = $elements.Cells.Item($r, $colRight)
    $leftValue = $leftCell.Value2
    $rightValue = $rightCell.Value2
    # Skip rows where both sides already hold the same (e.g. both blank)
    # value so untouched cells are left exactly as they were.
    if ($leftValue -ne $rightValue) {
        $leftCell.Value2 = $rightValue
        $rightCell.Value2 = $leftValue
    }
}

# --- 3. Swap the stored column widths to match the swapped content ----
# Column AK originally stored width 24.98046875 ("Mapping: RIM Mapping"),
# column AL originally stored width 85.5234375 ("Mapping: Spécification...").
# After the content swap, AK should carry the wider column and AL the
# narrower one. ColumnWidth is specified in characters and gets quantized
# to whole pixels by the host (stored_width = (Round(chars*6)+5)/6), so we
# pick the character width whose quantized result lands closest to the
# original target stored widths.
$elements.Columns.Item($colLeft).ColumnWidth = 84.66666666666667   # -> stored ~85.5 (was AL's width)
$elements.Columns.Item($colRight).ColumnWidth = 24.166666666666668 # -> stored ~25  (was AK's width)
